# Update Deepak Hooda (Kings XI Punjab) innings stats.
# Rows 2, 3, 5, 6 get new runs/balls/fours/sixes figures; row 4 is untouched.
# The source values are numeric-looking but stored as text in the workbook,
# so we force a text number format before writing so Excel doesn't silently
# convert them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cells = @("C2","D2","E2","F2","C3","D3","F3","C5","D5","C6","D6","E6")
foreach ($cell in $cells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Row 2: runs, balls, fours, sixes
$ws.Range("C2").Value = "23"
$ws.Range("D2").Value = "16"
$ws.Range("E2").Value = "1"
$ws.Range("F2").Value = "1"

# Row 3: runs, balls, sixes (fours unchanged)
$ws.Range("C3").Value = "15"
$ws.Range("D3").Value = "22"
$ws.Range("F3").Value = "0"

# Row 5: runs, balls (fours, sixes unchanged)
$ws.Range("C5").Value = "1"
$ws.Range("D5").Value = "1"

# Row 6: runs, balls, fours (sixes unchanged)
$ws.Range("C6").Value = "0"
$ws.Range("D6").Value = "2"
$ws.Range("E6").Value = "0"
